# Add three new GlobalConstant rows (TimeSecToGetOneEnergy, RequiredEnergyToPlay,
# RefillEnergyDiamond) to the "GlobalConstantIntTable" sheet and make that sheet
# the active/selected tab (moving the selection away from the
# "GlobalConstantFloatTable" sheet).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GlobalConstantIntTable")

# New shared strings must be minted in the same order as the target workbook:
# RequiredEnergyToPlay (50), TimeSecToGetOneEnergy (51), RefillEnergyDiamond (52).
$ws1.Range("A8").Value = "RequiredEnergyToPlay"
$ws1.Range("A7").Value = "TimeSecToGetOneEnergy"
$ws1.Range("A9").Value = "RefillEnergyDiamond"

$ws1.Range("C7").Value = 576
$ws1.Range("C8").Value = 5
$ws1.Range("C9").Value = 30

# Switch the active sheet/selection to GlobalConstantIntTable, cell C9.
$ws1.Activate()
$ws1.Range("C9").Select() | Out-Null
